$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the Wed Jul 26 22:41:05 UTC 2023 refresh.
# NumberFormat="@" + ClearFormats() ensures numeric-looking text (e.g. "1.002")
# is stored as text rather than being auto-converted to a number by Excel,
# while leaving the cell style/format untouched (matches original formatting).

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '29.473.77'
$cell.ClearFormats()
$cell = $ws.Range('E2')
$cell.NumberFormat = "@"
$cell.Value = '  +0.81%  '
$cell.ClearFormats()
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '1.869.95'
$cell.ClearFormats()
$cell = $ws.Range('E3')
$cell.NumberFormat = "@"
$cell.Value = '  +0.43%  '
$cell.ClearFormats()
$cell = $ws.Range('D4')
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.ClearFormats()
$cell = $ws.Range('E4')
$cell.NumberFormat = "@"
$cell.Value = '  +0.19%  '
$cell.ClearFormats()
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '0.7169'
$cell.ClearFormats()
$cell = $ws.Range('E5')
$cell.NumberFormat = "@"
$cell.Value = '  +0.56%  '
$cell.ClearFormats()
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '239.01'
$cell.ClearFormats()
$cell = $ws.Range('E6')
$cell.NumberFormat = "@"
$cell.Value = '  +0.40%  '
$cell.ClearFormats()
$cell = $ws.Range('E7')
$cell.NumberFormat = "@"
$cell.Value = '  +0.06%  '
$cell.ClearFormats()
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.07788'
$cell.ClearFormats()
$cell = $ws.Range('E8')
$cell.NumberFormat = "@"
$cell.Value = '  -4.99%  '
$cell.ClearFormats()
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.3069'
$cell.ClearFormats()
$cell = $ws.Range('E9')
$cell.NumberFormat = "@"
$cell.Value = '  +0.77%  '
$cell.ClearFormats()
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '25.22'
$cell.ClearFormats()
$cell = $ws.Range('E10')
$cell.NumberFormat = "@"
$cell.Value = '  +8.53%  '
$cell.ClearFormats()
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.08232'
$cell.ClearFormats()
$cell = $ws.Range('E11')
$cell.NumberFormat = "@"
$cell.Value = '  +0.71%  '
$cell.ClearFormats()
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '1.891.19'
$cell.ClearFormats()
$cell = $ws.Range('E12')
$cell.NumberFormat = "@"
$cell.Value = '  +0.94%  '
$cell.ClearFormats()
$cell = $ws.Range('B13')
$cell.NumberFormat = "@"
$cell.Value = 'Polkadot'
$cell.ClearFormats()
$cell = $ws.Range('C13')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell.ClearFormats()
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '5.225'
$cell.ClearFormats()
$cell = $ws.Range('E13')
$cell.NumberFormat = "@"
$cell.Value = '  +0.91%  '
$cell.ClearFormats()
$cell = $ws.Range('B14')
$cell.NumberFormat = "@"
$cell.Value = 'Polygon'
$cell.ClearFormats()
$cell = $ws.Range('C14')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell.ClearFormats()
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '0.7200'
$cell.ClearFormats()
$cell = $ws.Range('E14')
$cell.NumberFormat = "@"
$cell.Value = '  +1.41%  '
$cell.ClearFormats()
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '89.89'
$cell.ClearFormats()
$cell = $ws.Range('E15')
$cell.NumberFormat = "@"
$cell.Value = '  +0.19%  '
$cell.ClearFormats()
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '29.492.07'
$cell.ClearFormats()
$cell = $ws.Range('E16')
$cell.NumberFormat = "@"
$cell.Value = '  +0.75%  '
$cell.ClearFormats()
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '5.819'
$cell.ClearFormats()
$cell = $ws.Range('E17')
$cell.NumberFormat = "@"
$cell.Value = '  +0.38%  '
$cell.ClearFormats()
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '0.000007846'
$cell.ClearFormats()
$cell = $ws.Range('E18')
$cell.NumberFormat = "@"
$cell.Value = '  -1.17%  '
$cell.ClearFormats()
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '240.65'
$cell.ClearFormats()
$cell = $ws.Range('E19')
$cell.NumberFormat = "@"
$cell.Value = '  +1.36%  '
$cell.ClearFormats()
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '13.29'
$cell.ClearFormats()
$cell = $ws.Range('E20')
$cell.NumberFormat = "@"
$cell.Value = '  -0.70%  '
$cell.ClearFormats()
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '2.126.42'
$cell.ClearFormats()
$cell = $ws.Range('E21')
$cell.NumberFormat = "@"
$cell.Value = '  +0.84%  '
$cell.ClearFormats()
$cell = $ws.Range('E22')
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.ClearFormats()
$cell = $ws.Range('E23')
$cell.NumberFormat = "@"
$cell.Value = '  +0.16%  '
$cell.ClearFormats()
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '7.712'
$cell.ClearFormats()
$cell = $ws.Range('E24')
$cell.NumberFormat = "@"
$cell.Value = '  +3.70%  '
$cell.ClearFormats()
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '0.1561'
$cell.ClearFormats()
$cell = $ws.Range('E25')
$cell.NumberFormat = "@"
$cell.Value = '  +7.03%  '
$cell.ClearFormats()
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '162.61'
$cell.ClearFormats()
$cell = $ws.Range('E26')
$cell.NumberFormat = "@"
$cell.Value = '  -0.08%  '
$cell.ClearFormats()
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '8.952'
$cell.ClearFormats()
$cell = $ws.Range('E27')
$cell.NumberFormat = "@"
$cell.Value = '  -0.12%  '
$cell.ClearFormats()
$cell = $ws.Range('E28')
$cell.NumberFormat = "@"
$cell.Value = '  +1.10%  '
$cell.ClearFormats()
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '1.932'
$cell.ClearFormats()
$cell = $ws.Range('E29')
$cell.NumberFormat = "@"
$cell.Value = '  -1.54%  '
$cell.ClearFormats()
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '1.357'
$cell.ClearFormats()
$cell = $ws.Range('E30')
$cell.NumberFormat = "@"
$cell.Value = '  -4.78%  '
$cell.ClearFormats()
$cell = $ws.Range('E31')
$cell.NumberFormat = "@"
$cell.Value = '  -0.12%  '
$cell.ClearFormats()
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '4.326'
$cell.ClearFormats()
$cell = $ws.Range('E32')
$cell.NumberFormat = "@"
$cell.Value = '  -1.77%  '
$cell.ClearFormats()
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '4.072'
$cell.ClearFormats()
$cell = $ws.Range('E33')
$cell.NumberFormat = "@"
$cell.Value = '  +1.09%  '
$cell.ClearFormats()
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '0.05251'
$cell.ClearFormats()
$cell = $ws.Range('E34')
$cell.NumberFormat = "@"
$cell.Value = '  +0.45%  '
$cell.ClearFormats()
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '1.196'
$cell.ClearFormats()
$cell = $ws.Range('E35')
$cell.NumberFormat = "@"
$cell.Value = '  +2.18%  '
$cell.ClearFormats()
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.7153'
$cell.ClearFormats()
$cell = $ws.Range('E36')
$cell.NumberFormat = "@"
$cell.Value = '  +0.83%  '
$cell.ClearFormats()
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$cell = $ws.Range('E37')
$cell.NumberFormat = "@"
$cell.Value = '  -0.14%  '
$cell.ClearFormats()
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '2.673'
$cell.ClearFormats()
$cell = $ws.Range('E38')
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell.ClearFormats()
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.01867'
$cell.ClearFormats()
$cell = $ws.Range('E39')
$cell.NumberFormat = "@"
$cell.Value = '  +0.33%  '
$cell.ClearFormats()
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '2.721'
$cell.ClearFormats()
$cell = $ws.Range('E40')
$cell.NumberFormat = "@"
$cell.Value = '  -0.42%  '
$cell.ClearFormats()
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '1.173.42'
$cell.ClearFormats()
$cell = $ws.Range('E41')
$cell.NumberFormat = "@"
$cell.Value = '  +2.77%  '
$cell.ClearFormats()
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '0.9066'
$cell.ClearFormats()
$cell = $ws.Range('E42')
$cell.NumberFormat = "@"
$cell.Value = '  -1.93%  '
$cell.ClearFormats()
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '5.990'
$cell.ClearFormats()
$cell = $ws.Range('E43')
$cell.NumberFormat = "@"
$cell.Value = '  +1.43%  '
$cell.ClearFormats()
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '0.4304'
$cell.ClearFormats()
$cell = $ws.Range('E44')
$cell.NumberFormat = "@"
$cell.Value = '  +0.37%  '
$cell.ClearFormats()
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '71.44'
$cell.ClearFormats()
$cell = $ws.Range('E45')
$cell.NumberFormat = "@"
$cell.Value = '  +1.43%  '
$cell.ClearFormats()
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$cell = $ws.Range('E46')
$cell.NumberFormat = "@"
$cell.Value = '  +0.12%  '
$cell.ClearFormats()
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '102.23'
$cell.ClearFormats()
$cell = $ws.Range('E47')
$cell.NumberFormat = "@"
$cell.Value = '  -0.65%  '
$cell.ClearFormats()
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '0.5359'
$cell.ClearFormats()
$cell = $ws.Range('E48')
$cell.NumberFormat = "@"
$cell.Value = '  -0.68%  '
$cell.ClearFormats()
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '1.762'
$cell.ClearFormats()
$cell = $ws.Range('E49')
$cell.NumberFormat = "@"
$cell.Value = '  -0.88%  '
$cell.ClearFormats()
$cell = $ws.Range('E50')
$cell.NumberFormat = "@"
$cell.Value = '  -0.68%  '
$cell.ClearFormats()
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '7.010'
$cell.ClearFormats()
$cell = $ws.Range('E51')
$cell.NumberFormat = "@"
$cell.Value = '  +0.77%  '
$cell.ClearFormats()